$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '71.084.56'
    'E2' = '  -0.30%  '
    'D3' = '3.798.21'
    'E3' = '  -1.20%  '
    'E4' = '  -0.06%  '
    'D5' = '699.88'
    'E5' = '  -0.70%  '
    'D6' = '170.03'
    'E6' = '  -1.75%  '
    'D7' = '3.798.64'
    'E7' = '  -1.13%  '
    'E8' = '  -0.03%  '
    'E9' = '  -0.44%  '
    'E10' = '  -1.82%  '
    'E11' = '  +2.75%  '
    'E12' = '  +4.06%  '
    'E13' = '  -2.34%  '
    'D14' = '36.25'
    'E14' = '  -1.79%  '
    'D15' = '4.439.77'
    'E15' = '  -1.18%  '
    'D16' = '3.783.79'
    'E16' = '  -1.14%  '
    'D17' = '71.197.65'
    'E17' = '  -0.28%  '
    'E18' = '  -0.64%  '
    'D19' = '17.59'
    'E19' = '  +0.78%  '
    'E20' = '  +0.18%  '
    'D21' = '513.01'
    'E21' = '  +3.03%  '
    'D22' = '10.48'
    'E22' = '  -2.24%  '
    'D23' = '0.715'
    'E23' = '  -2.36%  '
    'D24' = '83.53'
    'E24' = '  -2.25%  '
    'E25' = '  -3.44%  '
    'D26' = '12.70'
    'E26' = '  +3.90%  '
    'D27' = '3.944.88'
    'E27' = '  -1.29%  '
    'D28' = '10.30'
    'E28' = '  -3.45%  '
    'E29' = '  +0.00%  '
    'E30' = '  -4.90%  '
    'D31' = '2.95'
    'E31' = '  -5.13%  '
    'E32' = '  +0.87%  '
    'D33' = '7.32'
    'E33' = '  -2.39%  '
    'D34' = '29.13'
    'E34' = '  -0.99%  '
    'E35' = '  -3.88%  '
    'E36' = '  +1.15%  '
    'D37' = '0.999'
    'E37' = '  -0.06%  '
    'D38' = '3.766.87'
    'E38' = '  -1.07%  '
    'D39' = '6.71'
    'E39' = '  +11.45%  '
    'E40' = '  -2.41%  '
    'E41' = '  +1.18%  '
    'E42' = '  -2.63%  '
    'E44' = '  -5.41%  '
    'E45' = '  -0.07%  '
    'D46' = '163.70'
    'E46' = '  -0.59%  '
    'B47' = 'OKB'
    'C47' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D47' = '49.13'
    'E47' = '  +0.13%  '
    'B48' = 'FLOKI'
    'C48' = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
    'D48' = '0.000304'
    'E48' = '  -4.30%  '
    'D49' = '425.04'
    'E49' = '  -1.14%  '
    'E50' = '  -0.99%  '
    'D51' = '1.38'
    'E51' = '  -0.41%  '
}

foreach ($ref in $updates.Keys) {
    $c = $ws.Range($ref)
    $savedStyle = $c.Style
    $c.NumberFormat = '@'
    $c.Value = $updates[$ref]
    $c.Style = $savedStyle
}
